$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")

# Copy formatting (style s="1") from row 822 down through the new rows
$ws.Range("A822:C822").Copy()
$ws.Range("A823:C856").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A823").Value = 'cs'
$ws.Range("B823").Value = 'market.atomizer.create.button'
$ws.Range("C823").Value = 'Nový atomizér'
$ws.Range("A824").Value = 'cs'
$ws.Range("B824").Value = 'market.atomizer.create.title'
$ws.Range("C824").Value = 'Nový atomizér'
$ws.Range("A825").Value = 'cs'
$ws.Range("B825").Value = 'shared.atomizer.create.name.label'
$ws.Range("C825").Value = 'Název atomizéru'
$ws.Range("A826").Value = 'cs'
$ws.Range("B826").Value = 'shared.atomizer.create.name.label.tooltip'
$ws.Range("C826").Value = 'Zadejte prosím název atomizéru ideálně stejně, jak jej pojmenoval výrobce. Ostatní uživatelé podle něj pak mohou (a budou) hledat.'
$ws.Range("A827").Value = 'cs'
$ws.Range("B827").Value = 'shared.atomizer.create.code.label'
$ws.Range("C827").Value = 'Kód atomizéru'
$ws.Range("A828").Value = 'cs'
$ws.Range("B828").Value = 'shared.atomizer.create.code.label.tooltip'
$ws.Range("C828").Value = 'Kód je dobrovolná hodnota určená hlavně pro importy/exporty a jednoznačenému rozeznání daného atomizéru.'
$ws.Range("A829").Value = 'cs'
$ws.Range("B829").Value = 'shared.atomizer.create.cost.label'
$ws.Range("C829").Value = 'Hodnota atomizéru'
$ws.Range("A830").Value = 'cs'
$ws.Range("B830").Value = 'shared.atomizer.create.cost.label.tooltip'
$ws.Range("C830").Value = 'Zadejte prosím hodnotu atomizéru na tržišti. Mějte prosím na paměti, že prodávané atomizéry by měly mít relativně podobné ceny.'
$ws.Range("A831").Value = 'cs'
$ws.Range("B831").Value = 'shared.atomizer.create.coilMin.label'
$ws.Range("C831").Value = 'Nejmenší spirálka'
$ws.Range("A832").Value = 'cs'
$ws.Range("B832").Value = 'shared.atomizer.create.coilMax.label'
$ws.Range("C832").Value = 'Největší spirálka'
$ws.Range("A833").Value = 'cs'
$ws.Range("B833").Value = 'shared.atomizer.create.coilMin.label.tooltip'
$ws.Range("C833").Value = 'Nejmenší doporučená velikost spriálky do atomizéru; pokud neuvádí výrobce, použijte zkušenost a zapište údaj dle svého názoru, ovšem bez extrémů.'
$ws.Range("A834").Value = 'cs'
$ws.Range("B834").Value = 'shared.atomizer.create.coilMax.label.tooltip'
$ws.Range("C834").Value = 'Největší doporučená velikost spriálky do atomizéru; pokud neuvádí výrobce, použijte zkušenost a zapište údaj dle svého názoru, ovšem bez extrémů.'
$ws.Range("A835").Value = 'cs'
$ws.Range("B835").Value = 'shared.atomizer.create.typeId.label'
$ws.Range("C835").Value = 'Typ atomizéru'
$ws.Range("A836").Value = 'cs'
$ws.Range("B836").Value = 'shared.atomizer.create.typeId.label.tooltip'
$ws.Range("C836").Value = 'Typ atomizéru by měl odpovídat jeho fyzické konstrukci; pokud např. výrobce udá atomizér jako RDA, ale reálně je to RDTA, uveďte prosím RDTA.'
$ws.Range("A837").Value = 'cs'
$ws.Range("B837").Value = 'shared.atomizer.create.drawIds.label'
$ws.Range("C837").Value = 'Typy potahů'
$ws.Range("A838").Value = 'cs'
$ws.Range("B838").Value = 'shared.atomizer.create.drawIds.label.tooltip'
$ws.Range("C838").Value = 'Zadejte prosím typy potahů, pro které je daný atomizér vhodný, ideálně z popisu výrobce, případně z vlastní praxe.'
$ws.Range("A839").Value = 'cs'
$ws.Range("B839").Value = 'shared.atomizer.create.squonk.label'
$ws.Range("C839").Value = 'Squonkovací atomizér?'
$ws.Range("A840").Value = 'cs'
$ws.Range("B840").Value = 'shared.atomizer.create.vendorId.label'
$ws.Range("C840").Value = 'Výrobce atomizéru'
$ws.Range("A841").Value = 'cs'
$ws.Range("B841").Value = 'shared.atomizer.create.create'
$ws.Range("C841").Value = 'Vytvořit atomizér'
$ws.Range("A842").Value = 'cs'
$ws.Range("B842").Value = 'shared.aroma.create.tasteIds.label'
$ws.Range("C842").Value = 'Příchutě'
$ws.Range("A843").Value = 'cs'
$ws.Range("B843").Value = 'shared.aroma.create.tasteIds.label.tooltip'
$ws.Range("C843").Value = 'Zadejte prosím příchutě v aromatu a to jak obecné (např. ovocné), tak i konkrétní (např. hruška); toto později umožní celkové hodnocení chutě z atomizéru, která bude založena na příchutích zde zadaných.'
$ws.Range("A844").Value = 'cs'
$ws.Range("B844").Value = 'common.taste.lemonade'
$ws.Range("C844").Value = 'Limonáda'
$ws.Range("A845").Value = 'cs'
$ws.Range("B845").Value = 'shared.atomizer.create.success'
$ws.Range("C845").Value = 'Atomizér byl úspěšně vytvořen.'
$ws.Range("A846").Value = 'cs'
$ws.Range("B846").Value = 'market.vendor.menu'
$ws.Range("C846").Value = 'Výrobci'
$ws.Range("A847").Value = 'cs'
$ws.Range("B847").Value = 'Vendor.list.total'
$ws.Range("C847").Value = 'Počet výrobců [{{data.total}}] ({{data.from}}-{{data.to}})'
$ws.Range("A848").Value = 'cs'
$ws.Range("B848").Value = 'market.vendor.index.title'
$ws.Range("C848").Value = 'Seznam výrobců'
$ws.Range("A849").Value = 'cs'
$ws.Range("B849").Value = 'common.filter.Vendor.filter.title'
$ws.Range("C849").Value = 'Vyhledat výrobce'
$ws.Range("A850").Value = 'cs'
$ws.Range("B850").Value = 'common.filter.Vendor.filter.id.label'
$ws.Range("C850").Value = 'Jméno výrobce'
$ws.Range("A851").Value = 'cs'
$ws.Range("B851").Value = 'shared.vendor.create.create'
$ws.Range("C851").Value = 'Vytvořit výrobce'
$ws.Range("A852").Value = 'cs'
$ws.Range("B852").Value = 'market.vendor.create.button'
$ws.Range("C852").Value = 'Vytvořit výrobce'
$ws.Range("A853").Value = 'cs'
$ws.Range("B853").Value = 'market.vendor.create.title'
$ws.Range("C853").Value = 'Nový výrobce'
$ws.Range("A854").Value = 'cs'
$ws.Range("B854").Value = 'shared.vendor.create.name.label'
$ws.Range("C854").Value = 'Zadejte jméno výrobce'
$ws.Range("A855").Value = 'cs'
$ws.Range("B855").Value = 'shared.vendor.create.name.label.tooltip'
$ws.Range("C855").Value = 'Zadejte prosím co nejpřesnějí jméno výrobce. Také se ujistěte, že už v systému není zadaný s podobným názevm.'
$ws.Range("A856").Value = 'cs'
$ws.Range("B856").Value = 'shared.vendor.create.success'
$ws.Range("C856").Value = 'Výrobce [{{name}}] byl úspěšně vytvořen.'

# Update view state to match: selection on B839, scrolled so row 823 is at top
$win = $excel.ActiveWindow
$win.ScrollRow = 823
$win.ScrollColumn = 1
$ws.Range("B839").Select()
